$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E")
$ws.Range("B3:H12").Value = 0
$ws.Range("B13:H27").Value = 1

$wsG = $wb.Worksheets.Item("G")
$wsG.Activate()
$excel.Goto($wsG.Range("A7"), $true)
$wsG.Range("E43").Select()

$ws.Activate()
$ws.Range("G17").Select()
